$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-30 02:25:20"

for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
